$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.25
$ws.Range("G2").Value = 3.55
$ws.Range("H2").Value = 2.8
$ws.Range("I2").Value = 3.05
$ws.Range("J2").Value = 2.68
$ws.Range("K2").Value = 2.92
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.18
$ws.Range("N2").Value = 2.1
$ws.Range("O2").Value = 1.8
$ws.Range("R2").Value = 1.11
$ws.Range("S2").Value = 7.6
$ws.Range("T2").Value = 2.42
$ws.Range("U2").Value = 1.56
$ws.Range("V2").Value = 1.48
$ws.Range("W2").Value = 1.39
$ws.Range("X2").Value = 6.4
$ws.Range("Y2").Value = 7.4
$ws.Range("Z2").Value = 18
$ws.Range("AA2").Value = 60
$ws.Range("AB2").Value = 8.199999999999999
$ws.Range("AC2").Value = 7.4
$ws.Range("AD2").Value = 16.5
$ws.Range("AE2").Value = 60
$ws.Range("AF2").Value = 22
$ws.Range("AG2").Value = 18.5
$ws.Range("AH2").Value = 34
$ws.Range("AI2").Value = 140
$ws.Range("AJ2").Value = 80
$ws.Range("AK2").Value = 75
$ws.Range("AL2").Value = 150
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 140
$ws.Range("AO2").Value = 100
$ws.Range("I3").Value = 1.72
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.03
$ws.Range("O3").Value = 1.11
$ws.Range("R3").Value = 1.74
$ws.Range("S3").Value = 1.83
$ws.Range("T3").Value = 1.01
$ws.Range("U3").Value = 2.4
$ws.Range("V3").Value = 2.38
$ws.Range("W3").Value = 1.17
$ws.Range("X3").Value = 55
$ws.Range("Y3").Value = 24
$ws.Range("Z3").Value = 21
$ws.Range("AA3").Value = 26
$ws.Range("AB3").Value = 48
$ws.Range("AC3").Value = 19
$ws.Range("AD3").Value = 17
$ws.Range("AE3").Value = 21
$ws.Range("AF3").Value = 75
$ws.Range("AG3").Value = 32
$ws.Range("AH3").Value = 25
$ws.Range("AI3").Value = 34
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 75
$ws.Range("AL3").Value = 65
$ws.Range("AM3").Value = 85
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000
$ws.Range("F6").Value = 1.86
$ws.Range("G6").Value = 2.02
$ws.Range("H6").Value = 4.1
$ws.Range("I6").Value = 7.6
$ws.Range("K6").Value = 3.5
$ws.Range("P7").Value = 1.05
$ws.Range("G9").Value = 2.04
$ws.Range("F10").Value = 2.92
$ws.Range("I10").Value = 2.58
$ws.Range("Q10").Value = 1.71
$ws.Range("F12").Value = 1.85
$ws.Range("I12").Value = 4.9
$ws.Range("P12").Value = 2.24
$ws.Range("F13").Value = 1.29
$ws.Range("G13").Value = 1.32
$ws.Range("Q13").Value = 1.68
$ws.Range("R13").Value = 1.51
$ws.Range("T13").Value = 2.28
$ws.Range("X13").Value = 28
$ws.Range("Y13").Value = 85
$ws.Range("AD13").Value = 190
$ws.Range("AH13").Value = 80
$ws.Range("AJ13").Value = 9.4
$ws.Range("AK13").Value = 18
$ws.Range("AL13").Value = 150
$ws.Range("AN13").Value = 5
$ws.Range("F14").Value = 1.43
$ws.Range("G14").Value = 1.46
$ws.Range("H14").Value = 8.800000000000001
$ws.Range("J14").Value = 4.8
$ws.Range("K14").Value = 5.2
$ws.Range("P14").Value = 2.14
$ws.Range("I15").Value = 1.85
$ws.Range("J15").Value = 3.9
$ws.Range("Q15").Value = 1.77
$ws.Range("G16").Value = 6.6
$ws.Range("AG17").Value = 11
$ws.Range("AK17").Value = 14.5
$ws.Range("F18").Value = 1.93
$ws.Range("G18").Value = 1.97
$ws.Range("I18").Value = 4.3
$ws.Range("P19").Value = 2.04
$ws.Range("Z19").Value = 100
$ws.Range("AD19").Value = 22
$ws.Range("I21").Value = 10.5
$ws.Range("J21").Value = 5.7
$ws.Range("K21").Value = 5.8
$ws.Range("Q21").Value = 1.63
$ws.Range("T21").Value = 1.99
$ws.Range("U21").Value = 1.92
$ws.Range("X21").Value = 27
$ws.Range("Y21").Value = 38
$ws.Range("AA21").Value = 440
$ws.Range("AC21").Value = 13
$ws.Range("AD21").Value = 70
$ws.Range("AL21").Value = 85
$ws.Range("P23").Value = 2.2
$ws.Range("Q23").Value = 1.73
$ws.Range("F24").Value = 2.02
$ws.Range("G24").Value = 2.08
$ws.Range("H24").Value = 3.65
$ws.Range("F25").Value = 1.35
$ws.Range("Q25").Value = 1.42
$ws.Range("L31").Value = 1.39
$ws.Range("G32").Value = 2.12
$ws.Range("H32").Value = 3.95
